$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 361
$ws1.Range("F4").Value = 1819
$ws1.Range("F10").Value = 3591
$ws1.Range("F11").Value = 145
$ws1.Range("F14").Value = 55
$ws1.Range("F15").Value = 61
$ws1.Range("F17").Value = 119
$ws1.Range("F18").Value = 788
$ws1.Range("F19").Value = 12
$ws1.Range("F20").Value = 214
$ws1.Range("F21").Value = 135
$ws1.Range("F23").Value = 71
$ws1.Range("F25").Value = 2848
$ws1.Range("F26").Value = 5306
$ws1.Range("F28").Value = 77
$ws1.Range("F30").Value = 3114
$ws1.Range("F31").Value = 303
$ws1.Range("F32").Value = 2302
$ws1.Range("F36").Value = 144
$ws1.Range("F37").Value = 196
$ws1.Range("F39").Value = 48
$ws1.Range("F40").Value = 475
$ws1.Range("F41").Value = 821
$ws1.Range("F42").Value = 33
$ws1.Range("F45").Value = 46
$ws1.Range("F46").Value = 506

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 361
$ws4.Range("F4").Value = 1819
$ws4.Range("F10").Value = 3591
$ws4.Range("F11").Value = 145
$ws4.Range("F15").Value = 55
$ws4.Range("F16").Value = 61
$ws4.Range("F18").Value = 119
$ws4.Range("F19").Value = 788
$ws4.Range("F20").Value = 12
$ws4.Range("F21").Value = 214
$ws4.Range("F22").Value = 135
$ws4.Range("F24").Value = 71
$ws4.Range("F26").Value = 2848
$ws4.Range("F27").Value = 5307
$ws4.Range("F29").Value = 77
$ws4.Range("F31").Value = 3114
$ws4.Range("F32").Value = 303
$ws4.Range("F33").Value = 2302
$ws4.Range("F37").Value = 144
$ws4.Range("F38").Value = 196
$ws4.Range("F40").Value = 48
$ws4.Range("F41").Value = 475
$ws4.Range("F42").Value = 821
$ws4.Range("F43").Value = 33
$ws4.Range("F46").Value = 46
$ws4.Range("F47").Value = 506
